$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.428.09'
$ws.Range("E2").Value = '  +0.65%  '
$ws.Range("D3").Value = '1.873.27'
$ws.Range("E3").Value = '  -0.59%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.021'
$ws.Range("E4").Value = '  +1.24%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.17'
$ws.Range("E5").Value = '  +0.45%  '
$ws.Range("E6").Value = '  +0.99%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5109'
$ws.Range("E7").Value = '  -0.52%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3950'
$ws.Range("E8").Value = '  +1.26%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08446'
$ws.Range("E9").Value = '  +0.82%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.110'
$ws.Range("E10").Value = '  -1.16%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '41.97'
$ws.Range("E11").Value = '  +0.59%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.250'
$ws.Range("E12").Value = '  +0.13%  '
$ws.Range("D13").Value = '1.880.82'
$ws.Range("E13").Value = '  +0.95%  '
$ws.Range("E14").Value = '  -0.56%  '
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.231'
$ws.Range("E15").Value = '  -0.42%  '
$ws.Range("B16").Value = 'BinanceUSD'
$ws.Range("C16").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.021'
$ws.Range("E16").Value = '  +1.20%  '
$ws.Range("E17").Value = '  +0.61%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '90.97'
$ws.Range("E18").Value = '  -0.09%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06772'
$ws.Range("E19").Value = '  +1.29%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.74'
$ws.Range("E20").Value = '  -0.39%  '
$ws.Range("E21").Value = '  +1.14%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.941'
$ws.Range("E22").Value = '  -1.52%  '
$ws.Range("D23").Value = '28.482.14'
$ws.Range("E23").Value = '  +0.67%  '
$ws.Range("E24").Value = '  +0.14%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.296'
$ws.Range("E25").Value = '  +0.65%  '
$ws.Range("D26").Value = '2.088.82'
$ws.Range("E26").Value = '  -0.25%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '161.89'
$ws.Range("E27").Value = '  +0.89%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.75'
$ws.Range("E28").Value = '  +0.37%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.348'
$ws.Range("E29").Value = '  -4.66%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '127.07'
$ws.Range("E30").Value = '  +0.37%  '
$ws.Range("E31").Value = '  -0.66%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.037'
$ws.Range("E32").Value = '  -0.24%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.767'
$ws.Range("E33").Value = '  -1.86%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.634'
$ws.Range("E34").Value = '  +0.65%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02435'
$ws.Range("E35").Value = '  -0.31%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06463'
$ws.Range("E36").Value = '  -1.62%  '
$ws.Range("E37").Value = '  -1.58%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.840'
$ws.Range("E38").Value = '  -6.76%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.260'
$ws.Range("E39").Value = '  +1.30%  '
$ws.Range("E40").Value = '  -1.45%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6368'
$ws.Range("E41").Value = '  -1.97%  '
$ws.Range("E42").Value = '  -0.65%  '
$ws.Range("E43").Value = '  -0.24%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6028'
$ws.Range("E44").Value = '  -1.19%  '
$ws.Range("E45").Value = '  -0.70%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.714'
$ws.Range("E46").Value = '  +0.48%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.990'
$ws.Range("E47").Value = '  -1.26%  '
$ws.Range("E48").Value = '  -5.80%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.205'
$ws.Range("E49").Value = '  -2.61%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '121.84'
$ws.Range("E50").Value = '  +0.52%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06837'
$ws.Range("E51").Value = '  -1.04%  '
